# add save column in s_vals sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1, "sum") onto the
# new header cell H1, then give it its own text ("Save").
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Save column values for rows 2-11
$saveValues = @(0, 0, 1, 0, 0, 0, 0, 1, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
